$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 74.609651
$ws.Range("H2").Value = 223.828953
$ws.Range("I2").Value = 0.1061386348809139
$ws.Range("J2").Value = 0.1061386348809139
$ws.Range("M2").Value = 0.242595
$ws.Range("N2").Value = 0.727785
$ws.Range("O2").Value = 0.03230700759563258
$ws.Range("P2").Value = 0.03230700759563257
$ws.Range("Q2").Value = 18.099928284345
$ws.Range("R2").Value = 162.899354559105
$ws.Range("S2").Value = 0.003429021683287759
$ws.Range("T2").Value = 0.003429021683287758
$ws.Range("G3").Value = 74.609651
$ws.Range("H3").Value = 223.828953
$ws.Range("I3").Value = 0.1061386348809139
$ws.Range("J3").Value = 0.1061386348809139
$ws.Range("O3").Value = 0.4893229309549773
$ws.Range("P3").Value = 0.4893229309549771
$ws.Range("Q3").Value = 274.1420706313847
$ws.Range("R3").Value = 2467.278635682462
$ws.Range("S3").Value = 0.05193606790748898
$ws.Range("T3").Value = 0.05193606790748898
$ws.Range("G4").Value = 74.609651
$ws.Range("H4").Value = 223.828953
$ws.Range("I4").Value = 0.1061386348809139
$ws.Range("J4").Value = 0.1061386348809139
$ws.Range("M4").Value = 3.574634666666667
$ws.Range("N4").Value = 10.723904
$ws.Range("O4").Value = 0.4760434029044768
$ws.Range("P4").Value = 0.4760434029044767
$ws.Range("Q4").Value = 266.7022449325013
$ws.Range("R4").Value = 2400.320204392513
$ws.Range("S4").Value = 0.05052659692834605
$ws.Range("T4").Value = 0.05052659692834605
$ws.Range("G5").Value = 74.609651
$ws.Range("H5").Value = 223.828953
$ws.Range("I5").Value = 0.1061386348809139
$ws.Range("J5").Value = 0.1061386348809139
$ws.Range("M5").Value = 0.017471
$ws.Range("N5").Value = 0.052413
$ws.Range("O5").Value = 0.002326658544913526
$ws.Range("P5").Value = 0.002326658544913525
$ws.Range("Q5").Value = 1.303505212621
$ws.Range("R5").Value = 11.731546913589
$ws.Range("S5").Value = 0.0002469483617911352
$ws.Range("T5").Value = 0.0002469483617911352
$ws.Range("G6").Value = 597.374756
$ws.Range("I6").Value = 0.8498168837991085
$ws.Range("J6").Value = 0.8498168837991086
$ws.Range("M6").Value = 0.242595
$ws.Range("N6").Value = 0.727785
$ws.Range("O6").Value = 0.03230700759563258
$ws.Range("P6").Value = 0.03230700759563257
$ws.Range("Q6").Value = 144.92012893182
$ws.Range("R6").Value = 1304.28116038638
$ws.Range("S6").Value = 0.02745504051979461
$ws.Range("T6").Value = 0.0274550405197946
$ws.Range("G7").Value = 597.374756
$ws.Range("I7").Value = 0.8498168837991085
$ws.Range("J7").Value = 0.8498168837991086
$ws.Range("O7").Value = 0.4893229309549773
$ws.Range("P7").Value = 0.4893229309549771
$ws.Range("Q7").Value = 2194.964731208275
$ws.Range("S7").Value = 0.4158348883556051
$ws.Range("T7").Value = 0.415834888355605
$ws.Range("G8").Value = 597.374756
$ws.Range("I8").Value = 0.8498168837991085
$ws.Range("J8").Value = 0.8498168837991086
$ws.Range("M8").Value = 3.574634666666667
$ws.Range("N8").Value = 10.723904
$ws.Range("O8").Value = 0.4760434029044768
$ws.Range("P8").Value = 0.4760434029044767
$ws.Range("Q8").Value = 2135.396511789142
$ws.Range("R8").Value = 19218.56860610228
$ws.Range("S8").Value = 0.4045497212094059
$ws.Range("T8").Value = 0.4045497212094059
$ws.Range("G9").Value = 597.374756
$ws.Range("I9").Value = 0.8498168837991085
$ws.Range("J9").Value = 0.8498168837991086
$ws.Range("M9").Value = 0.017471
$ws.Range("N9").Value = 0.052413
$ws.Range("O9").Value = 0.002326658544913526
$ws.Range("P9").Value = 0.002326658544913525
$ws.Range("Q9").Value = 10.436734362076
$ws.Range("R9").Value = 93.930609258684
$ws.Range("S9").Value = 0.001977233714302981
$ws.Range("T9").Value = 0.001977233714302981
$ws.Range("G10").Value = 30.48438
$ws.Range("H10").Value = 91.45313999999999
$ws.Range("I10").Value = 0.04336664808137267
$ws.Range("J10").Value = 0.04336664808137267
$ws.Range("M10").Value = 0.242595
$ws.Range("N10").Value = 0.727785
$ws.Range("O10").Value = 0.03230700759563258
$ws.Range("P10").Value = 0.03230700759563257
$ws.Range("Q10").Value = 7.395358166099999
$ws.Range("R10").Value = 66.55822349489999
$ws.Range("S10").Value = 0.001401046628962032
$ws.Range("T10").Value = 0.001401046628962031
$ws.Range("G11").Value = 30.48438
$ws.Range("H11").Value = 91.45313999999999
$ws.Range("I11").Value = 0.04336664808137267
$ws.Range("J11").Value = 0.04336664808137267
$ws.Range("O11").Value = 0.4893229309549773
$ws.Range("P11").Value = 0.4893229309549771
$ws.Range("Q11").Value = 112.01032229884
$ws.Range("R11").Value = 1008.09290068956
$ws.Range("S11").Value = 0.02122029534487031
$ws.Range("T11").Value = 0.02122029534487031
$ws.Range("G12").Value = 30.48438
$ws.Range("H12").Value = 91.45313999999999
$ws.Range("I12").Value = 0.04336664808137267
$ws.Range("J12").Value = 0.04336664808137267
$ws.Range("M12").Value = 3.574634666666667
$ws.Range("N12").Value = 10.723904
$ws.Range("O12").Value = 0.4760434029044768
$ws.Range("P12").Value = 0.4760434029044767
$ws.Range("Q12").Value = 108.97052153984
$ws.Range("R12").Value = 980.73469385856
$ws.Range("S12").Value = 0.02064440672521754
$ws.Range("T12").Value = 0.02064440672521754
$ws.Range("G13").Value = 30.48438
$ws.Range("H13").Value = 91.45313999999999
$ws.Range("I13").Value = 0.04336664808137267
$ws.Range("J13").Value = 0.04336664808137267
$ws.Range("M13").Value = 0.017471
$ws.Range("N13").Value = 0.052413
$ws.Range("O13").Value = 0.002326658544913526
$ws.Range("P13").Value = 0.002326658544913525
$ws.Range("Q13").Value = 0.53259260298
$ws.Range("R13").Value = 4.793333426819999
$ws.Range("S13").Value = 0.0001008993823227835
$ws.Range("T13").Value = 0.0001008993823227835
$ws.Range("E14").Value = 3
$ws.Range("F14").Value = 1
$ws.Range("G14").Value = 0.4764796666666666
$ws.Range("H14").Value = 1.429439
$ws.Range("I14").Value = 0.0006778332386049212
$ws.Range("J14").Value = 0.0006778332386049213
$ws.Range("M14").Value = 0.242595
$ws.Range("N14").Value = 0.727785
$ws.Range("O14").Value = 0.03230700759563258
$ws.Range("P14").Value = 0.03230700759563257
$ws.Range("Q14").Value = 0.115591584735
$ws.Range("R14").Value = 1.040324262615
$ws.Range("S14").Value = 0.00002189876358818142
$ws.Range("T14").Value = 0.00002189876358818142
$ws.Range("E15").Value = 3
$ws.Range("F15").Value = 1
$ws.Range("G15").Value = 0.4764796666666666
$ws.Range("H15").Value = 1.429439
$ws.Range("I15").Value = 0.0006778332386049212
$ws.Range("J15").Value = 0.0006778332386049213
$ws.Range("O15").Value = 0.4893229309549773
$ws.Range("P15").Value = 0.4893229309549771
$ws.Range("Q15").Value = 1.750753698522889
$ws.Range("R15").Value = 15.756783286706
$ws.Range("S15").Value = 0.0003316793470128644
$ws.Range("T15").Value = 0.0003316793470128644
$ws.Range("E16").Value = 3
$ws.Range("F16").Value = 1
$ws.Range("G16").Value = 0.4764796666666666
$ws.Range("H16").Value = 1.429439
$ws.Range("I16").Value = 0.0006778332386049212
$ws.Range("J16").Value = 0.0006778332386049213
$ws.Range("M16").Value = 3.574634666666667
$ws.Range("N16").Value = 10.723904
$ws.Range("O16").Value = 0.4760434029044768
$ws.Range("P16").Value = 0.4760434029044767
$ws.Range("Q16").Value = 1.703240734428444
$ws.Range("R16").Value = 15.329166609856
$ws.Range("S16").Value = 0.0003226780415072488
$ws.Range("T16").Value = 0.0003226780415072488
$ws.Range("E17").Value = 3
$ws.Range("F17").Value = 1
$ws.Range("G17").Value = 0.4764796666666666
$ws.Range("H17").Value = 1.429439
$ws.Range("I17").Value = 0.0006778332386049212
$ws.Range("J17").Value = 0.0006778332386049213
$ws.Range("M17").Value = 0.017471
$ws.Range("N17").Value = 0.052413
$ws.Range("O17").Value = 0.002326658544913526
$ws.Range("P17").Value = 0.002326658544913525
$ws.Range("Q17").Value = 0.008324576256333333
$ws.Range("R17").Value = 0.074921186307
$ws.Range("S17").Value = 0.000001577086496626549
$ws.Range("T17").Value = 0.000001577086496626549
